# Applies the "new fiscal year column" update to the RAVN financials sheet:
# inserts a new column D (period ending 2019-01-31) and shifts the existing
# D:K data one column to the right (E:L), preserving per-cell number formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new blank column at D; existing D:K shift to E:L.
$ws.Columns("D").Insert()

# 2) The freshly inserted column D copied its formatting from column C
#    (label column). Re-copy the number formats from column E (the old
#    column D, now shifted one to the right) so the new column D matches
#    the other data columns.
$src = $ws.Range("E5:E102")
$dst = $ws.Range("D5:D102")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Populate the new column D with the new fiscal year's figures.
$ws.Range("D7").Value = 43496
$ws.Range("D8").Value = 406700
$ws.Range("D9").Value = 274100
$ws.Range("D10").Value = 132500
$ws.Range("D12").Value = 26200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 351500
$ws.Range("D18").Value = 55100
$ws.Range("D20").Value = 6400
$ws.Range("D21").Value = 76700
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 61600
$ws.Range("D24").Value = 9700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 51900
$ws.Range("D27").Value = 51800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -6400
$ws.Range("D33").Value = 51800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 51800

$ws.Range("D38").Value = 43496
$ws.Range("D41").Value = 65800
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 55500
$ws.Range("D44").Value = 54100
$ws.Range("D45").Value = 7700
$ws.Range("D46").Value = 183100
$ws.Range("D47").Value = 300
$ws.Range("D48").Value = 106600
$ws.Range("D49").Value = 67200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 360200
$ws.Range("D57").Value = 8300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 24800
$ws.Range("D60").Value = 33100
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 18200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 51300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 286000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 309000
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43496
$ws.Range("D81").Value = 51800
$ws.Range("D83").Value = 15100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 66000
$ws.Range("D91").Value = -14100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -16400
$ws.Range("D96").Value = -18800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -23800
$ws.Range("D101").Value = -500
$ws.Range("D102").Value = 25300
